# Tenth day first commit: append rows 96-101 to the TaskList sheet,
# matching the formatting/banding of the preceding rows (94:95).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TaskList")

# New task rows to append, in order (row, day-serial, task text, duration text)
$newRows = @(
  @{ A = 95;  C = "Quiz";                                      F = "60Minutes" },
  @{ A = 96;  C = "Create RegisterHandler class";               F = "20Minutes" },
  @{ A = 97;  C = "Created exception.jsp page";                 F = "10Minutes" },
  @{ A = 98;  C = "Testing Webflow Registration page";          F = "120Minutes" },
  @{ A = 99;  C = "Created Address DAO and DTO and tested";     F = "60Minutes" },
  @{ A = 100; C = "Configured Billing page with webflow";       F = "60Minutes" }
)

$startRow = 96
$taskDate = Get-Date -Year 2017 -Month 3 -Day 8 -Hour 0 -Minute 0 -Second 0

for ($i = 0; $i -lt $newRows.Count; $i++) {
  $destRow = $startRow + $i

  # Bring over the same look (borders/fonts/number formats/wrap) as the two
  # rows immediately above the table end, preserving the odd/even banding:
  # even destination rows mirror row 94's formatting, odd ones mirror row 95's.
  if ($destRow % 2 -eq 0) {
    $srcRow = 94
  } else {
    $srcRow = 95
  }
  $ws.Range("A$srcRow`:H$srcRow").Copy()
  $ws.Range("A$destRow`:H$destRow").PasteSpecial(-4122)

  $row = $newRows[$i]
  $ws.Cells.Item($destRow, 1).Value = $row.A
  $ws.Cells.Item($destRow, 2).Value = $taskDate
  $ws.Cells.Item($destRow, 3).Value = $row.C
  $ws.Cells.Item($destRow, 6).Value = $row.F
}

$excel.CutCopyMode = $false

# Match the author's scroll/selection position from the commit.
$ws.Application.ActiveWindow.ScrollRow = 88
$ws.Range("D107").Select()
